# [GEN MCU] Sample 제작
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GEN MCU")

# Make "GEN MCU" the active sheet (was "GEN Main" before).
$ws.Activate() | Out-Null

# Insert a new "Version" column before the existing "Category" column (old column D).
$ws.Columns("D:D").Insert() | Out-Null

# New column header.
$ws.Cells.Item(4, 4).Value = "Version"

# Fill the Version column for the existing rows with "V1.0".
$ws.Cells.Item(5, 4).Value = "V1.0"
$ws.Cells.Item(6, 4).Value = "V1.0"
$ws.Cells.Item(7, 4).Value = "V1.0"
$ws.Cells.Item(8, 4).Value = "V1.0"
$ws.Cells.Item(9, 4).Value = "V1.0"
$ws.Cells.Item(10, 4).Value = "V1.0"
$ws.Cells.Item(11, 4).Value = "V1.0"
$ws.Cells.Item(12, 4).Value = "V1.0"

# New issue rows (9 and 10 -> sheet rows 13 and 14).
$ws.Cells.Item(13, 3).Value = 43213
$ws.Cells.Item(13, 4).Value = "V2.0"
$ws.Cells.Item(13, 6).Value = "'+5V CON 변경 - 기존 version과 동일한 CON 사용?"

$ws.Cells.Item(14, 3).Value = 43213
$ws.Cells.Item(14, 4).Value = "V2.0"
$ws.Cells.Item(14, 6).Value = "R, C PCB PAD 확장 - PAD 간격은 유지?"

$ws.Cells.Item(13, 5).Value = "고민중"
$ws.Cells.Item(14, 5).Value = "고민중"

# Update selection to match the saved view.
$ws.Range("E14").Select() | Out-Null
